# Generate Report for Handback
# Adds a new handback row (73c57d2d-a97f-4404-8339-9002bdb9b00a) to the
# Overview / zh-cn / de-de report sheets, extending each table by one row.

$wb = $excel.ActiveWorkbook

$fileId   = "73c57d2d-a97f-4404-8339-9002bdb9b00a"
$mdName   = "$fileId.md"
$mdPath   = "e2e\$fileId.md"
$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Cells.Item(4, 1).Value = $mdName
$wsOverview.Cells.Item(4, 2).Value = $mdPath
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(4, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/05a6fa204a1911e7c0c51f7912c6c3394d084950/e2e/$mdName", "", "", $mdPath) | Out-Null
$wsOverview.Cells.Item(4, 3).Value = ".md"
$wsOverview.Cells.Item(4, 5).Value = $statusInSync
$wsOverview.Cells.Item(4, 6).Value = $statusInSync
$wsOverview.Cells.Item(4, 7).Value = "2016-09-07 11:02:43"
$wsOverview.Cells.Item(4, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Shared per-locale row values (columns D..P are identical in shape
# between zh-cn and de-de, only G/H/J/K differ)
# ---------------------------------------------------------------------
function Add-LocaleRow($ws, $xlfName, $handoffDate, $handbackDate, $hyperlinkRepoSuffix, $hyperlinkSha) {

    $lo = $ws.ListObjects.Item(1)
    $lo.ListRows.Add() | Out-Null

    $ws.Cells.Item(4, 1).Value = $mdName
    $ws.Cells.Item(4, 2).Value = ".md"
    $ws.Cells.Item(4, 3).Value = $statusInSync
    $ws.Cells.Item(4, 4).Value = "e2e"
    $ws.Cells.Item(4, 5).Value = "ht"
    $ws.Cells.Item(4, 6).Value = "True"
    $ws.Cells.Item(4, 7).Value = $xlfName
    $ws.Cells.Item(4, 8).Value = $handoffDate
    $ws.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Cells.Item(4, 9).Value = $mdName
    $ws.Cells.Item(4, 10).Value = $xlfName
    $ws.Cells.Item(4, 11).Value = $handbackDate
    $ws.Cells.Item(4, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Cells.Item(4, 12).Value = ""
    $ws.Cells.Item(4, 13).Value = "True"
    $ws.Cells.Item(4, 14).Value = ""
    $ws.Cells.Item(4, 15).Value = "False"
    $ws.Cells.Item(4, 16).Value = ""

    $ws.Hyperlinks.Add($ws.Cells.Item(4, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/05a6fa204a1911e7c0c51f7912c6c3394d084950/e2e/$mdName", "", "", $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(4, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-$hyperlinkRepoSuffix/blob/$hyperlinkSha/e2e/$mdName", "", "", $mdName) | Out-Null
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Add-LocaleRow $wsZhCn "$fileId.05a6fa204a1911e7c0c51f7912c6c3394d084950.zh-cn.xlf" "2016-09-07 11:02:37" "2016-09-07 11:02:57" "zhcn" "2e5abff8b4cd4be239151d9de7bbece41e33c485"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
Add-LocaleRow $wsDeDe "$fileId.05a6fa204a1911e7c0c51f7912c6c3394d084950.de-de.xlf" "2016-09-07 11:02:43" "2016-09-07 11:03:16" "dede" "b5ffbe4f985a95066c977b9e3b5922f8c40555d5"
